$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122 ; used to copy the "section label" font (fontId 3,
# cellXfs index 4 in the original file) from a cell that already carries it
# (A16) onto the C-column cells that become "秦森" below, without
# introducing brand-new font/style table entries.
$xlPasteFormats = -4122

# Row 8: 用户组流转 -> 执行人 秦森 (50% complete)
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C8").Value = "秦森"
$ws.Range("D8").Value = 0.5
$ws.Range("D8").NumberFormat = "0%"
$excel.CutCopyMode = $false

# Row 9: 消息开始事件学习及应用 -> 执行人 秦森 (progress stays 100%)
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C9").Value = "秦森"
$excel.CutCopyMode = $false

# Row 10: 错误结束事件学习及应用 -> 执行人 秦森 (now 100% complete)
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Value = "秦森"
$ws.Range("D10").Value = 1
$ws.Range("D10").NumberFormat = "0%"
$excel.CutCopyMode = $false

# Row 11: 排他网关、并行网关学习及应用 (执行人 unchanged) now 100% complete
$ws.Range("D11").Value = 1
$ws.Range("D11").NumberFormat = "0%"

# Row 16: 子流程学习及流程应用 -> 执行人 秦森
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C16").Value = "秦森"
$excel.CutCopyMode = $false

# Leave the cursor where the author ended up after this edit.
$ws.Range("E10").Select() | Out-Null
